$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header D1 from "ALUOP" to "ALUOP[1:0]"
$ws.Range("D1").Value = "ALUOP[1:0]"

# Set column D width (matches width="11.5" custom width added in diff; the
# engine snaps stored width to 1/7-character increments, so 10.785714... is
# the ColumnWidth input that lands closest on the target stored width)
$ws.Columns.Item(4).ColumnWidth = 10.785714285714286

# Add new row 13 data: SLT instruction with ALUOP value 2
$ws.Range("A13").Value = "SLT"
$ws.Range("B13").Value = '`1010'
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0

# Update selection to D13 as in the diff
$ws.Range("D13").Select() | Out-Null
